# Apply the Color-wise data sheet update:
#  - remove the rows for cars that were dropped from this export
#    (zeekr-001, maxus-mifa-7, bmw-x2, NIO EL6)
#  - rename the remaining "Car Names" entries to their simplified form
#  - narrow column A from 48 to 29 characters

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows bottom-to-top so earlier row numbers stay valid.
$ws.Rows(16).Delete()   # NIO EL6-2024-1054.xlsx
$ws.Rows(8).Delete()    # bmw-x2-2022-1065.xlsx
$ws.Rows(4).Delete()    # maxus-mifa-7-2024-1060.xlsx
$ws.Rows(3).Delete()    # zeekr-001-2024-1037.xlsx

# Rename the remaining car names to their simplified labels.
$ws.Range("A2").Value = "ford-tourneo-custom- 2024"
$ws.Range("A3").Value = "ford-tourneo-custom- 2024"
$ws.Range("A4").Value = "vw-passat- 2024"
$ws.Range("A5").Value = "skoda-kodiaq- 2024"
$ws.Range("A6").Value = "renault-rafale-hev- 2022"
$ws.Range("A7").Value = "mercedes-benz-e-class- 2024"
$ws.Range("A8").Value = "suzuki-swift- 2024"
$ws.Range("A9").Value = "dacia-duster- 2024"
$ws.Range("A10").Value = "renault-espace- 2022"
$ws.Range("A11").Value = "Toyota C-HR- 2024"
$ws.Range("A12").Value = "Honda CR-V- 2024"
$ws.Range("A13").Value = "Honda CR-V- 2024"

# Column A width 48 -> 29 (character units, snaps to raw width 29).
$ws.Columns(1).ColumnWidth = 28.14
